# Generate Report for Handback
#
# - zh-cn / de-de sheets: copy the "Source File Name" (A) and
#   "Latest Handoff File" (D) hyperlinked values into the
#   "Latest Target File" (F) and "Latest Handback File" (G) columns
#   for each data row, since the handed-back files are now in sync
#   with the handed-off ones.
# - Update the "Status" column (shared across the Overview summary
#   sheet too) from "Ready for handoff" to "Handed back: in sync with en-US".
# - Stamp "Latest Handback DateTime" (H) for both language sheets -
#   zh-cn finished slightly before de-de.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$zhHandbackTime = "2016-03-20 20:51:19"
$deHandbackTime = "2016-03-20 20:51:24"

function Copy-HyperlinkedValue {
    param($ws, [string]$srcCell, [string]$dstCell)

    $srcRange = $ws.Range($srcCell)
    $dstRange = $ws.Range($dstCell)

    $srcHyperlink = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $srcRange.Address()) {
            $srcHyperlink = $hl
        }
    }

    $dstRange.Value = $srcRange.Value

    if ($srcHyperlink -ne $null) {
        $ws.Hyperlinks.Add($dstRange, $srcHyperlink.Address, "", "", $srcHyperlink.TextToDisplay) | Out-Null
    }
}

# ---- Overview sheet: roll the new status text into the summary grid ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

Copy-HyperlinkedValue $wsZh "A2" "F2"
Copy-HyperlinkedValue $wsZh "D2" "G2"
Copy-HyperlinkedValue $wsZh "A3" "F3"
Copy-HyperlinkedValue $wsZh "D3" "G3"

$wsZh.Range("H2").Value = $zhHandbackTime
$wsZh.Range("H3").Value = $zhHandbackTime

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

Copy-HyperlinkedValue $wsDe "A2" "F2"
Copy-HyperlinkedValue $wsDe "D2" "G2"
Copy-HyperlinkedValue $wsDe "A3" "F3"
Copy-HyperlinkedValue $wsDe "D3" "G3"

$wsDe.Range("H2").Value = $deHandbackTime
$wsDe.Range("H3").Value = $deHandbackTime
